# "add update to week after rollercoaster week"
# Append six new weigh-in readings (rows 269-274) to the raw_data sheet,
# right after the existing data which ended at row 268.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("raw_data")

# New rows of data: date/time serial, time-of-day serial, weight (kg), AM/PM
$newRows = @(
    @{ Row = 269; DateTime = 44118.327777777777; Time = 0.32777777777777778;  Weight = 71.3 },
    @{ Row = 270; DateTime = 44117.92291666667;  Time = 0.92291666666666661;  Weight = 71.7 },
    @{ Row = 271; DateTime = 44117.431944444441; Time = 0.43194444444444446; Weight = 71.2 },
    @{ Row = 272; DateTime = 44117.431250000001; Time = 0.43124999999999997; Weight = 71.2 },
    @{ Row = 273; DateTime = 44117.431250000001; Time = 0.43124999999999997; Weight = 71.900000000000006 },
    @{ Row = 274; DateTime = 44117.334027777775; Time = 0.33402777777777781; Weight = 71.900000000000006 }
)

foreach ($item in $newRows) {
    $r = $item.Row

    $ws.Range("A$r").Value = $item.DateTime
    $ws.Range("A$r").NumberFormat = "m/d/yy h:mm"

    $ws.Range("B$r").Value = $item.Time
    $ws.Range("B$r").NumberFormat = "h:mm"

    $ws.Range("C$r").Value = $item.Weight
}

# Column D: extend the AM/PM helper formula down through the new rows in one
# shot so the engine groups them into a single shared-formula block.
$ws.Range("D269:D274").Formula = "=IF(B269<TIME(12,0,0), ""AM"", ""PM"")"

# Move the selection/scroll position to just past the new data, matching
# where Excel leaves the cursor after typing the last row of a manual entry.
[void]$ws.Range("A275").Select()
